# Macroferia Regional de Talca - Melon (Hortaliza) - weekly update
# Inserts a new week's worth of price data (6 rows) at the top of the
# date-ordered block that starts at row 426, pushing the existing rows
# (426:449) down to (432:455). The new rows carry date 44931 (the newest
# reporting week) together with fresh Calameno/Tuna x Extra/Primera/Segunda
# price observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows right before the current row 426, shifting the rest
# of the table (426:449) down to (432:455).
$ws.Range("A426:A431").EntireRow.Insert()

$newRows = @(
    @{ Row = 426; Variedad = "Calameño"; Calidad = "Extra";   Volumen = 4000; Precio = 1000 },
    @{ Row = 427; Variedad = "Calameño"; Calidad = "Primera"; Volumen = 4000; Precio = 800  },
    @{ Row = 428; Variedad = "Calameño"; Calidad = "Segunda"; Volumen = 3000; Precio = 600  },
    @{ Row = 429; Variedad = "Tuna";     Calidad = "Extra";   Volumen = 3000; Precio = 1300 },
    @{ Row = 430; Variedad = "Tuna";     Calidad = "Primera"; Volumen = 3000; Precio = 1000 },
    @{ Row = 431; Variedad = "Tuna";     Calidad = "Segunda"; Volumen = 2000; Precio = 700  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = 44931
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = 100112027
    $ws.Cells.Item($row, 7).Value = "Melón"
    $ws.Cells.Item($row, 8).Value = $r.Variedad
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.Precio
    $ws.Cells.Item($row, 12).Value = $r.Precio
    $ws.Cells.Item($row, 13).Value = $r.Precio
    $ws.Cells.Item($row, 14).Value = "$/unidad"
    $ws.Cells.Item($row, 15).Value = "Región del Maule"
    $ws.Cells.Item($row, 16).Value = $r.Precio
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
